$d = $word.ActiveDocument
$d.Content.Find.Execute("Nucleo-F401RE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nucleo-L552ZE-Q", 2)
